$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each data row, shift all values one column to the left (B<-C, C<-D, ...)
# and clear the now-redundant last column, matching the ifoCAST full series update.
$cols = @("B","C","D","E","F","G","H","I","J","K")

# Row 2: populated through column K
$ws.Range("B2").Value2 = $ws.Range("C2").Value2
$ws.Range("C2").Value2 = $ws.Range("D2").Value2
$ws.Range("D2").Value2 = $ws.Range("E2").Value2
$ws.Range("E2").Value2 = $ws.Range("F2").Value2
$ws.Range("F2").Value2 = $ws.Range("G2").Value2
$ws.Range("G2").Value2 = $ws.Range("H2").Value2
$ws.Range("H2").Value2 = $ws.Range("I2").Value2
$ws.Range("I2").Value2 = $ws.Range("J2").Value2
$ws.Range("J2").Value2 = $ws.Range("K2").Value2
$ws.Range("K2").ClearContents()

# Row 3: populated through column K
$ws.Range("B3").Value2 = $ws.Range("C3").Value2
$ws.Range("C3").Value2 = $ws.Range("D3").Value2
$ws.Range("D3").Value2 = $ws.Range("E3").Value2
$ws.Range("E3").Value2 = $ws.Range("F3").Value2
$ws.Range("F3").Value2 = $ws.Range("G3").Value2
$ws.Range("G3").Value2 = $ws.Range("H3").Value2
$ws.Range("H3").Value2 = $ws.Range("I3").Value2
$ws.Range("I3").Value2 = $ws.Range("J3").Value2
$ws.Range("J3").Value2 = $ws.Range("K3").Value2
$ws.Range("K3").ClearContents()

# Row 4: populated through column K
$ws.Range("B4").Value2 = $ws.Range("C4").Value2
$ws.Range("C4").Value2 = $ws.Range("D4").Value2
$ws.Range("D4").Value2 = $ws.Range("E4").Value2
$ws.Range("E4").Value2 = $ws.Range("F4").Value2
$ws.Range("F4").Value2 = $ws.Range("G4").Value2
$ws.Range("G4").Value2 = $ws.Range("H4").Value2
$ws.Range("H4").Value2 = $ws.Range("I4").Value2
$ws.Range("I4").Value2 = $ws.Range("J4").Value2
$ws.Range("J4").Value2 = $ws.Range("K4").Value2
$ws.Range("K4").ClearContents()

# Row 5: populated through column K
$ws.Range("B5").Value2 = $ws.Range("C5").Value2
$ws.Range("C5").Value2 = $ws.Range("D5").Value2
$ws.Range("D5").Value2 = $ws.Range("E5").Value2
$ws.Range("E5").Value2 = $ws.Range("F5").Value2
$ws.Range("F5").Value2 = $ws.Range("G5").Value2
$ws.Range("G5").Value2 = $ws.Range("H5").Value2
$ws.Range("H5").Value2 = $ws.Range("I5").Value2
$ws.Range("I5").Value2 = $ws.Range("J5").Value2
$ws.Range("J5").Value2 = $ws.Range("K5").Value2
$ws.Range("K5").ClearContents()

# Row 6: populated through column K
$ws.Range("B6").Value2 = $ws.Range("C6").Value2
$ws.Range("C6").Value2 = $ws.Range("D6").Value2
$ws.Range("D6").Value2 = $ws.Range("E6").Value2
$ws.Range("E6").Value2 = $ws.Range("F6").Value2
$ws.Range("F6").Value2 = $ws.Range("G6").Value2
$ws.Range("G6").Value2 = $ws.Range("H6").Value2
$ws.Range("H6").Value2 = $ws.Range("I6").Value2
$ws.Range("I6").Value2 = $ws.Range("J6").Value2
$ws.Range("J6").Value2 = $ws.Range("K6").Value2
$ws.Range("K6").ClearContents()

# Row 7: populated through column K
$ws.Range("B7").Value2 = $ws.Range("C7").Value2
$ws.Range("C7").Value2 = $ws.Range("D7").Value2
$ws.Range("D7").Value2 = $ws.Range("E7").Value2
$ws.Range("E7").Value2 = $ws.Range("F7").Value2
$ws.Range("F7").Value2 = $ws.Range("G7").Value2
$ws.Range("G7").Value2 = $ws.Range("H7").Value2
$ws.Range("H7").Value2 = $ws.Range("I7").Value2
$ws.Range("I7").Value2 = $ws.Range("J7").Value2
$ws.Range("J7").Value2 = $ws.Range("K7").Value2
$ws.Range("K7").ClearContents()

# Row 8: populated through column K
$ws.Range("B8").Value2 = $ws.Range("C8").Value2
$ws.Range("C8").Value2 = $ws.Range("D8").Value2
$ws.Range("D8").Value2 = $ws.Range("E8").Value2
$ws.Range("E8").Value2 = $ws.Range("F8").Value2
$ws.Range("F8").Value2 = $ws.Range("G8").Value2
$ws.Range("G8").Value2 = $ws.Range("H8").Value2
$ws.Range("H8").Value2 = $ws.Range("I8").Value2
$ws.Range("I8").Value2 = $ws.Range("J8").Value2
$ws.Range("J8").Value2 = $ws.Range("K8").Value2
$ws.Range("K8").ClearContents()

# Row 9: populated through column K
$ws.Range("B9").Value2 = $ws.Range("C9").Value2
$ws.Range("C9").Value2 = $ws.Range("D9").Value2
$ws.Range("D9").Value2 = $ws.Range("E9").Value2
$ws.Range("E9").Value2 = $ws.Range("F9").Value2
$ws.Range("F9").Value2 = $ws.Range("G9").Value2
$ws.Range("G9").Value2 = $ws.Range("H9").Value2
$ws.Range("H9").Value2 = $ws.Range("I9").Value2
$ws.Range("I9").Value2 = $ws.Range("J9").Value2
$ws.Range("J9").Value2 = $ws.Range("K9").Value2
$ws.Range("K9").ClearContents()

# Row 10: populated through column K
$ws.Range("B10").Value2 = $ws.Range("C10").Value2
$ws.Range("C10").Value2 = $ws.Range("D10").Value2
$ws.Range("D10").Value2 = $ws.Range("E10").Value2
$ws.Range("E10").Value2 = $ws.Range("F10").Value2
$ws.Range("F10").Value2 = $ws.Range("G10").Value2
$ws.Range("G10").Value2 = $ws.Range("H10").Value2
$ws.Range("H10").Value2 = $ws.Range("I10").Value2
$ws.Range("I10").Value2 = $ws.Range("J10").Value2
$ws.Range("J10").Value2 = $ws.Range("K10").Value2
$ws.Range("K10").ClearContents()

# Row 11: populated through column K
$ws.Range("B11").Value2 = $ws.Range("C11").Value2
$ws.Range("C11").Value2 = $ws.Range("D11").Value2
$ws.Range("D11").Value2 = $ws.Range("E11").Value2
$ws.Range("E11").Value2 = $ws.Range("F11").Value2
$ws.Range("F11").Value2 = $ws.Range("G11").Value2
$ws.Range("G11").Value2 = $ws.Range("H11").Value2
$ws.Range("H11").Value2 = $ws.Range("I11").Value2
$ws.Range("I11").Value2 = $ws.Range("J11").Value2
$ws.Range("J11").Value2 = $ws.Range("K11").Value2
$ws.Range("K11").ClearContents()

# Row 12: populated through column K
$ws.Range("B12").Value2 = $ws.Range("C12").Value2
$ws.Range("C12").Value2 = $ws.Range("D12").Value2
$ws.Range("D12").Value2 = $ws.Range("E12").Value2
$ws.Range("E12").Value2 = $ws.Range("F12").Value2
$ws.Range("F12").Value2 = $ws.Range("G12").Value2
$ws.Range("G12").Value2 = $ws.Range("H12").Value2
$ws.Range("H12").Value2 = $ws.Range("I12").Value2
$ws.Range("I12").Value2 = $ws.Range("J12").Value2
$ws.Range("J12").Value2 = $ws.Range("K12").Value2
$ws.Range("K12").ClearContents()

# Row 13: populated through column K
$ws.Range("B13").Value2 = $ws.Range("C13").Value2
$ws.Range("C13").Value2 = $ws.Range("D13").Value2
$ws.Range("D13").Value2 = $ws.Range("E13").Value2
$ws.Range("E13").Value2 = $ws.Range("F13").Value2
$ws.Range("F13").Value2 = $ws.Range("G13").Value2
$ws.Range("G13").Value2 = $ws.Range("H13").Value2
$ws.Range("H13").Value2 = $ws.Range("I13").Value2
$ws.Range("I13").Value2 = $ws.Range("J13").Value2
$ws.Range("J13").Value2 = $ws.Range("K13").Value2
$ws.Range("K13").ClearContents()

# Row 14: populated through column K
$ws.Range("B14").Value2 = $ws.Range("C14").Value2
$ws.Range("C14").Value2 = $ws.Range("D14").Value2
$ws.Range("D14").Value2 = $ws.Range("E14").Value2
$ws.Range("E14").Value2 = $ws.Range("F14").Value2
$ws.Range("F14").Value2 = $ws.Range("G14").Value2
$ws.Range("G14").Value2 = $ws.Range("H14").Value2
$ws.Range("H14").Value2 = $ws.Range("I14").Value2
$ws.Range("I14").Value2 = $ws.Range("J14").Value2
$ws.Range("J14").Value2 = $ws.Range("K14").Value2
$ws.Range("K14").ClearContents()

# Row 15: populated through column K
$ws.Range("B15").Value2 = $ws.Range("C15").Value2
$ws.Range("C15").Value2 = $ws.Range("D15").Value2
$ws.Range("D15").Value2 = $ws.Range("E15").Value2
$ws.Range("E15").Value2 = $ws.Range("F15").Value2
$ws.Range("F15").Value2 = $ws.Range("G15").Value2
$ws.Range("G15").Value2 = $ws.Range("H15").Value2
$ws.Range("H15").Value2 = $ws.Range("I15").Value2
$ws.Range("I15").Value2 = $ws.Range("J15").Value2
$ws.Range("J15").Value2 = $ws.Range("K15").Value2
$ws.Range("K15").ClearContents()

# Row 16: populated through column J
$ws.Range("B16").Value2 = $ws.Range("C16").Value2
$ws.Range("C16").Value2 = $ws.Range("D16").Value2
$ws.Range("D16").Value2 = $ws.Range("E16").Value2
$ws.Range("E16").Value2 = $ws.Range("F16").Value2
$ws.Range("F16").Value2 = $ws.Range("G16").Value2
$ws.Range("G16").Value2 = $ws.Range("H16").Value2
$ws.Range("H16").Value2 = $ws.Range("I16").Value2
$ws.Range("I16").Value2 = $ws.Range("J16").Value2
$ws.Range("J16").ClearContents()

# Row 17: populated through column I
$ws.Range("B17").Value2 = $ws.Range("C17").Value2
$ws.Range("C17").Value2 = $ws.Range("D17").Value2
$ws.Range("D17").Value2 = $ws.Range("E17").Value2
$ws.Range("E17").Value2 = $ws.Range("F17").Value2
$ws.Range("F17").Value2 = $ws.Range("G17").Value2
$ws.Range("G17").Value2 = $ws.Range("H17").Value2
$ws.Range("H17").Value2 = $ws.Range("I17").Value2
$ws.Range("I17").ClearContents()

# Row 18: populated through column H
$ws.Range("B18").Value2 = $ws.Range("C18").Value2
$ws.Range("C18").Value2 = $ws.Range("D18").Value2
$ws.Range("D18").Value2 = $ws.Range("E18").Value2
$ws.Range("E18").Value2 = $ws.Range("F18").Value2
$ws.Range("F18").Value2 = $ws.Range("G18").Value2
$ws.Range("G18").Value2 = $ws.Range("H18").Value2
$ws.Range("H18").ClearContents()

# Row 19: populated through column G
$ws.Range("B19").Value2 = $ws.Range("C19").Value2
$ws.Range("C19").Value2 = $ws.Range("D19").Value2
$ws.Range("D19").Value2 = $ws.Range("E19").Value2
$ws.Range("E19").Value2 = $ws.Range("F19").Value2
$ws.Range("F19").Value2 = $ws.Range("G19").Value2
$ws.Range("G19").ClearContents()

# Row 20: populated through column F
$ws.Range("B20").Value2 = $ws.Range("C20").Value2
$ws.Range("C20").Value2 = $ws.Range("D20").Value2
$ws.Range("D20").Value2 = $ws.Range("E20").Value2
$ws.Range("E20").Value2 = $ws.Range("F20").Value2
$ws.Range("F20").ClearContents()

# Row 21: populated through column E
$ws.Range("B21").Value2 = $ws.Range("C21").Value2
$ws.Range("C21").Value2 = $ws.Range("D21").Value2
$ws.Range("D21").Value2 = $ws.Range("E21").Value2
$ws.Range("E21").ClearContents()

# Row 22: populated through column D
$ws.Range("B22").Value2 = $ws.Range("C22").Value2
$ws.Range("C22").Value2 = $ws.Range("D22").Value2
$ws.Range("D22").ClearContents()

# Row 23: populated through column C
$ws.Range("B23").Value2 = $ws.Range("C23").Value2
$ws.Range("C23").ClearContents()

# Row 24: populated through column B
$ws.Range("B24").ClearContents()

